$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = "'" + '68.329.57'
$ws.Range('E2').Value2 = "'" + '  -1.68%  '
$ws.Range('D3').Value2 = "'" + '3.828.44'
$ws.Range('E3').Value2 = "'" + '  -1.67%  '
$ws.Range('E4').Value2 = "'" + '  -0.38%  '
$ws.Range('D5').Value2 = "'" + '601.24'
$ws.Range('E5').Value2 = "'" + '  -0.59%  '
$ws.Range('D6').Value2 = "'" + '170.22'
$ws.Range('E6').Value2 = "'" + '  +0.89%  '
$ws.Range('D7').Value2 = "'" + '3.829.81'
$ws.Range('E7').Value2 = "'" + '  -1.84%  '
$ws.Range('E8').Value2 = "'" + '  -0.25%  '
$ws.Range('E9').Value2 = "'" + '  -1.25%  '
$ws.Range('E10').Value2 = "'" + '  -2.21%  '
$ws.Range('D11').Value2 = "'" + '6.48'
$ws.Range('E11').Value2 = "'" + '  +1.39%  '
$ws.Range('D12').Value2 = "'" + '0.457'
$ws.Range('E12').Value2 = "'" + '  -1.87%  '
$ws.Range('D13').Value2 = "'" + '0.0000263'
$ws.Range('E13').Value2 = "'" + '  +3.48%  '
$ws.Range('D14').Value2 = "'" + '37.10'
$ws.Range('E14').Value2 = "'" + '  -2.54%  '
$ws.Range('D15').Value2 = "'" + '4.473.66'
$ws.Range('E15').Value2 = "'" + '  -2.14%  '
$ws.Range('D16').Value2 = "'" + '3.830.09'
$ws.Range('E16').Value2 = "'" + '  -1.93%  '
$ws.Range('D17').Value2 = "'" + '68.314.38'
$ws.Range('E17').Value2 = "'" + '  -2.07%  '
$ws.Range('D18').Value2 = "'" + '18.50'
$ws.Range('E18').Value2 = "'" + '  -1.12%  '
$ws.Range('E19').Value2 = "'" + '  -2.39%  '
$ws.Range('E20').Value2 = "'" + '  -0.65%  '
$ws.Range('D21').Value2 = "'" + '11.12'
$ws.Range('E21').Value2 = "'" + '  -0.01%  '
$ws.Range('D22').Value2 = "'" + '469.21'
$ws.Range('E22').Value2 = "'" + '  -4.33%  '
$ws.Range('E23').Value2 = "'" + '  -1.15%  '
$ws.Range('E24').Value2 = "'" + '  -3.79%  '
$ws.Range('D25').Value2 = "'" + '83.11'
$ws.Range('E25').Value2 = "'" + '  -2.73%  '
$ws.Range('E26').Value2 = "'" + '  -2.21%  '
$ws.Range('E27').Value2 = "'" + '  -0.95%  '
$ws.Range('D28').Value2 = "'" + '10.01'
$ws.Range('E28').Value2 = "'" + '  -1.23%  '
$ws.Range('E29').Value2 = "'" + '  +0.06%  '
$ws.Range('E30').Value2 = "'" + '  -0.58%  '
$ws.Range('D31').Value2 = "'" + '3.979.76'
$ws.Range('E31').Value2 = "'" + '  -2.18%  '
$ws.Range('D32').Value2 = "'" + '7.70'
$ws.Range('E32').Value2 = "'" + '  -1.13%  '
$ws.Range('D33').Value2 = "'" + '31.56'
$ws.Range('E33').Value2 = "'" + '  -1.17%  '
$ws.Range('E34').Value2 = "'" + '  -4.45%  '
$ws.Range('D35').Value2 = "'" + '9.44'
$ws.Range('E35').Value2 = "'" + '  -0.78%  '
$ws.Range('D36').Value2 = "'" + '3.794.73'
$ws.Range('E36').Value2 = "'" + '  -2.09%  '
$ws.Range('D37').Value2 = "'" + '0.105'
$ws.Range('E37').Value2 = "'" + '  -2.27%  '
$ws.Range('D38').Value2 = "'" + '3.69'
$ws.Range('E38').Value2 = "'" + '  +12.64%  '
$ws.Range('E39').Value2 = "'" + '  -1.86%  '
$ws.Range('E40').Value2 = "'" + '  -1.21%  '
$ws.Range('E41').Value2 = "'" + '  -2.53%  '
$ws.Range('D42').Value2 = "'" + '0.999'
$ws.Range('E42').Value2 = "'" + '  -0.54%  '
$ws.Range('E43').Value2 = "'" + '  -3.81%  '
$ws.Range('E44').Value2 = "'" + '  -5.65%  '
$ws.Range('D45').Value2 = "'" + '8.78'
$ws.Range('E45').Value2 = "'" + '  +1.48%  '
$ws.Range('D46').Value2 = "'" + '0.000295'
$ws.Range('E46').Value2 = "'" + '  +9.67%  '
$ws.Range('B47').Value2 = 'Bittensor'
$ws.Range('C47').Value2 = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D47').Value2 = "'" + '417.71'
$ws.Range('E47').Value2 = "'" + '  -4.10%  '
$ws.Range('B48').Value2 = 'USDe'
$ws.Range('C48').Value2 = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D48').Value2 = "'" + '1.00'
$ws.Range('E48').Value2 = "'" + '  -0.01%  '
$ws.Range('D49').Value2 = "'" + '47.11'
$ws.Range('E49').Value2 = "'" + '  -1.97%  '
$ws.Range('D50').Value2 = "'" + '26.15'
$ws.Range('E50').Value2 = "'" + '  +3.28%  '
$ws.Range('D51').Value2 = "'" + '141.51'
$ws.Range('E51').Value2 = "'" + '  -1.53%  '
